$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-TextValue "D2" "305.44"
Set-TextValue "E2" "0.95%"
Set-TextValue "D3" "36.11"
Set-TextValue "E3" "-3.05%"
Set-TextValue "D4" "5.100"
Set-TextValue "E4" "2.03%"
Set-TextValue "D5" "0.07878"
Set-TextValue "E5" "0.73%"
Set-TextValue "D6" "2.130"
Set-TextValue "E6" "-3.72%"
Set-TextValue "D7" "7.928"
Set-TextValue "E7" "-1.27%"
Set-TextValue "D8" "0.9230"
Set-TextValue "E8" "0.87%"
Set-TextValue "D9" "0.09723"
Set-TextValue "E9" "-0.15%"
Set-TextValue "E10" "-1.41%"
Set-TextValue "D11" "0.08614"
Set-TextValue "E11" "-0.96%"
Set-TextValue "D12" "0.03559"
Set-TextValue "E12" "1.12%"
Set-TextValue "D13" "0.09923"
Set-TextValue "E13" "-0.40%"
Set-TextValue "D14" "0.001429"
Set-TextValue "E14" "-3.33%"
Set-TextValue "D15" "0.005622"
Set-TextValue "E15" "-0.57%"
Set-TextValue "D16" "3.464"
Set-TextValue "E16" "0.14%"
Set-TextValue "D17" "4.106"
Set-TextValue "E17" "1.86%"
Set-TextValue "D18" "2.635"
Set-TextValue "E18" "16.50%"
Set-TextValue "E19" "-1.89%"
Set-TextValue "D20" "0.1317"
Set-TextValue "E20" "1.31%"
Set-TextValue "D21" "5.170"
Set-TextValue "E21" "8.55%"
Set-TextValue "D22" "0.2203"
Set-TextValue "E22" "-4.03%"
Set-TextValue "D23" "0.04550"
Set-TextValue "E23" "-1.48%"
Set-TextValue "D24" "0.005054"
Set-TextValue "E24" "5.51%"
Set-TextValue "D25" "0.001233"
Set-TextValue "E25" "0.26%"
Set-TextValue "D27" "0.0004743"
Set-TextValue "D39" "0.01846"
Set-TextValue "E39" "4.39%"
Set-TextValue "D40" "0.04769"
Set-TextValue "E40" "0.37%"
Set-TextValue "D41" "0.007545"
Set-TextValue "E41" "-6.30%"
Set-TextValue "E42" "0.59%"
Set-TextValue "D43" "0.007724"
Set-TextValue "E43" "0.86%"
Set-TextValue "D44" "0.002228"
Set-TextValue "E44" "3.14%"
Set-TextValue "D45" "0.01104"
Set-TextValue "E45" "11.86%"
Set-TextValue "D46" "0.00006319"
Set-TextValue "E46" "5.23%"
Set-TextValue "E47" "-0.15%"
Set-TextValue "D48" "0.0005792"
Set-TextValue "E48" "-0.15%"
Set-TextValue "D49" "47.59"
Set-TextValue "E49" "501.89%"
Set-TextValue "E51" "-0.15%"
